$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.588.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.726.84'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +3.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("E6").Value = '  +6.59%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.725.43'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.03%  '
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("E11").Value = '  +4.46%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.157'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.34'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("E14").Value = '  +2.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.223.23'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.576.89'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.699.44'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.94'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '374.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.07%  '
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("E22").Value = '  +2.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.03%  '
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.34'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.58%  '
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '589.47'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.81%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.31'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("E33").Value = '  +3.74%  '
$ws.Range("E34").Value = '  +6.00%  '
$ws.Range("E35").Value = '  +2.88%  '
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.04'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("E40").Value = '  +2.76%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.47'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.90'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  -3.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.97'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("E48").Value = '  +5.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '155.44'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("E51").Value = '  +5.88%  '

Write-Host "Applied all changes"